$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains text-formatted numbers (e.g. "43.174.07" or "327.64").
# Force text number-format on D2:D51 before assigning so purely-numeric-looking
# strings (e.g. "327.64") are not auto-converted to numeric values by Excel,
# then restore the original style so no stray formatting is introduced.
$dRange = $ws.Range("D2:D51")
$dOrigStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.174.07"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "2.391.40"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "327.64"
$ws.Range("E5").Value = "  +10.88%  "
$ws.Range("E6").Value = "  -5.75%  "
$ws.Range("D7").Value = "0.648"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.649"
$ws.Range("E9").Value = "  +6.58%  "
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").Value = "0.0937"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "8.67"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "17.05"
$ws.Range("E14").Value = "  +11.74%  "
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "2.753.84"
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").Value = "2.386.59"
$ws.Range("E17").Value = "  +6.13%  "
$ws.Range("D18").Value = "43.182.16"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  +7.36%  "
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "76.47"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "3.74"
$ws.Range("E22").Value = "  +7.57%  "
$ws.Range("D23").Value = "271.36"
$ws.Range("E23").Value = "  +5.74%  "
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("D26").Value = "11.78"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "23.00"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "175.77"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.19"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "37.53"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D33").Value = "0.0929"
$ws.Range("E33").Value = "  +4.46%  "
$ws.Range("D34").Value = "5.93"
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("D36").Value = "4.93"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("D37").Value = "4.13"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").Value = "0.0366"
$ws.Range("E38").Value = "  -2.95%  "
$ws.Range("D39").Value = "0.109"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +16.54%  "
$ws.Range("E41").Value = "  +18.96%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "69.81"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").Value = "121.65"
$ws.Range("E44").Value = "  +13.87%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "12.38"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "89.12"
$ws.Range("E47").Value = "  +45.08%  "
$ws.Range("D48").Value = "9.36"
$ws.Range("E48").Value = "  +7.93%  "
$ws.Range("D49").Value = "5.54"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").Value = "0.488"
$ws.Range("E51").Value = "  +11.25%  "

# Restore original (default) style/number-format for column D
$dRange.Style = $dOrigStyle

